$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.598.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.866.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4788"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3813"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07357"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9349"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07815"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.865.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.448"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.575"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008845"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.694.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.26%  "
$ws.Range("E21").Value = "  +1.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.106"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.027"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.950"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08890"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.328"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("E32").Value = "  +2.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7612"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.606"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.704"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.135"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.83%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02037"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.07%  "
$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5693"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05387"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.982"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.057"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.556"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1528"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4902"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "105.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.49%  "
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.665"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "67.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9115"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.92%  "
